$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- 1. Update the first three rows: 100 -> 0M, 0.01 -> 0M, 229 -> 0M ---
$t.Rows(1).Cells(1).Range.Text = "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"

# --- 2. Insert 10 new rows right after row 3 (before the row that currently
#        holds "0", i.e. the former row 4), in order:
#        146, 0.00004, 0.00027, 0.00009, 0.00005, 0.00008, 0.00008,
#        0.00027, 0.01105, 100.0
#        Rows.Add(beforeRow) always inserts immediately before the anchor,
#        so pushing the same anchor repeatedly reverses the order; insert
#        the values back-to-front so the final layout reads forward. ---
$newRowTexts = @("146", "0.00004", "0.00027", "0.00009", "0.00005", "0.00008", "0.00008", "0.00027", "0.01105", "100.0")
$anchorRow = $t.Rows(4)
for ($i = $newRowTexts.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells(1).Range.Text = $newRowTexts[$i]
}

# --- 3. Collapse the final three multi-run rows (each a tab-separated list
#        of values) down to just their first value. These rows are now at
#        index 44, 45, 46 after the 10-row insertion above. ---
$lastIndex = $t.Rows.Count
$t.Rows($lastIndex - 2).Cells(1).Range.Text = "100"
$t.Rows($lastIndex - 1).Cells(1).Range.Text = "0.01"
$t.Rows($lastIndex).Cells(1).Range.Text = "229"
